$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.258.14'
$ws.Range("E2").Value = '  -2.02%  '

# Row 3
$ws.Range("D3").Value = '1.878.60'
$ws.Range("E3").Value = '  -1.47%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.04%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.06%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4853'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.04%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2879'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.84%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06586'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.54%  '

# Row 10
$ws.Range("D10").Value = '1.887.78'
$ws.Range("E10").Value = '  -1.11%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.73'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.86%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07292'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.03%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.163'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.84%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.95%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6554'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.17%  '

# Row 16
$ws.Range("D16").Value = '30.230.65'
$ws.Range("E16").Value = '  -2.04%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.34'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.07%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007733'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.66%  '

# Row 20
$ws.Range("D20").Value = '2.127.64'
$ws.Range("E20").Value = '  -1.24%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.296'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.79%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9992'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.10%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '193.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.47%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.123'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.24%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.276'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.81%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.09'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.35%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.914'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.67%  '

# Row 29
$ws.Range("E29").Value = '  +0.56%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.264'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.43%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09119'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.55%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.039'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.27%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05093'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.59%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7179'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.56%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.096'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.67%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.698'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.44%  '

# Row 37
$ws.Range("E37").Value = '  -2.49%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.638'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.09%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9196'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.72%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.040'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.95%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '106.32'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.03%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4271'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.84%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.784'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.98%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.60%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.393'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.75%  '

# Row 46
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1316'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.18%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.53%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.922'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.11%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05752'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.60%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.83'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.36%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3813'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.21%  '
